$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# Set the "Execute" column (A) to "N" for every scenario row except row 12
# (cancelVisit), which stays "Y" so only that test runs after the new
# OpsAPI-based visit cleanup is introduced.
for ($r = 2; $r -le 14; $r++) {
    if ($r -ne 12) {
        $ws.Cells.Item($r, 1).Value = "N"
    }
}

# Reflect the author's last selected cell when saving the workbook.
$ws.Range("A12").Select()
